$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 16) mirroring the existing rows' layout.
$row = 16

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 42625.886678240742

$ws.Cells.Item($row, 2).Value = -22
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 48
$ws.Cells.Item($row, 5).Value = 40
$ws.Cells.Item($row, 6).Value = 60
$ws.Cells.Item($row, 7).Value = 13423
$ws.Cells.Item($row, 8).Value = 9531
$ws.Cells.Item($row, 9).Value = 424
$ws.Cells.Item($row, 10).Value = 97
$ws.Cells.Item($row, 11).Value = 93
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 3
$ws.Cells.Item($row, 14).Value = "Named"
